$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"
$ws.Range("G7").ClearContents() | Out-Null
$ws.Range("G8").ClearContents() | Out-Null

$ws.Range("G4").Select() | Out-Null
